$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Features  To dos")
$ws.Activate()

$ws.Range("A12").Value = "Fatigue moment logger"
$ws.Range("C12").Value = "Thomas"
$ws.Range("D12").Value = "Added new way to log data during main"
$ws.Range("F12").Value = "Thomas_workspace"

$ws.Range("B12").Value = 44665
$ws.Range("E12").Value = 44665
$ws.Range("E12").NumberFormat = "d-mmm"

$ws.Range("F13").Select()
